$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the measurement columns (new data set: Measurement3 / Measurement4)
$ws.Range("C1").Value = "Measurement3"
$ws.Range("D1").Value = "Measurement4"

# Updated measurement values for the existing Cell_type/Replicate rows
$data = @(
    @(1.5, 1.5),
    @(1.2, 1.4),
    @(2.6, 2.4),
    @(2.5, 2.7),
    @(1.9, 1.4),
    @(2.2000000000000002, 2.4),
    @(2.2000000000000002, 2.5),
    @(1.4, 1.6),
    @(1.5, 1.1000000000000001),
    @(1.9, 1.4),
    @(2.2000000000000002, 2.5),
    @(2.2000000000000002, 2.8),
    @(2.8, 2.7),
    @(2.9, 2.2000000000000002),
    @(3, 2.8),
    @(1.7, 1.5),
    @(2.6, 2.4),
    @(2.6, 2.2000000000000002),
    @(2.2000000000000002, 2.5),
    @(2, 2.6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][1]
}
